# Weekly update: insert a new weekly price record for
# "Macroferia Regional de Talca - Choclo" as the new row 163,
# pushing the existing rows 163-176 down to 164-177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 163 (shifts rows 163:176 down to 164:177,
# mirroring the existing row's formatting, e.g. the date style on column D).
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with this week's data.
$ws.Cells.Item(163, 1).Value  = 5
$ws.Cells.Item(163, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(163, 3).Value  = "Maule"
$ws.Cells.Item(163, 4).Value  = 44585
$ws.Cells.Item(163, 5).Value  = 7
$ws.Cells.Item(163, 6).Value  = 100112024
$ws.Cells.Item(163, 7).Value  = "Choclo"
$ws.Cells.Item(163, 8).Value  = "Choclero"
$ws.Cells.Item(163, 9).Value  = "Primera"
$ws.Cells.Item(163, 10).Value = 60000
$ws.Cells.Item(163, 11).Value = 120
$ws.Cells.Item(163, 12).Value = 140
$ws.Cells.Item(163, 13).Value = 130
$ws.Cells.Item(163, 14).Value = "$/unidad"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 130
$ws.Cells.Item(163, 17).Value = 1
$ws.Cells.Item(163, 18).Value = "Hortaliza"
